# Rename the header cell from "Categoría de negocio" to "lista_negocio".
# (The rest of the apparent diff - xmlns additions, theme renames, row
# height/column width tweaks, selection, etc. - are just artifacts of the
# file being re-saved by a newer Excel build and carry no content change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "lista_negocio"
